$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows continuing the Ticker list (2024-02-08 data update)
$ws.Range("A294").Value = "IMX-USD"
$ws.Range("A295").Value = "TAO-USD"
$ws.Range("A296").Value = "MNT-USD"
